$d = $word.ActiveDocument

# The document ends with (after the last "Requisito" line, LOQ4095):
#   [empty paragraph]
#   "Ver no Jupiter Salvar em pdf Salvar em docx"
#   [empty paragraph]
#   [empty paragraph with a page break before it]
#   [empty paragraph]                      <- stays
#   [empty paragraph with a page break before it]  <- stays
#
# The edit removes the 4-paragraph block consisting of the empty paragraph
# right before the "Ver no Jupiter..." line, the "Ver no Jupiter..." line
# itself, the empty paragraph right after it, and the following empty
# page-break paragraph. The two final (empty) paragraphs that originally
# closed the document are left untouched.

# Find the paragraph holding the unique anchor text.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -gt 0) {
    $startParaIndex = $anchorIndex - 1
    $endParaIndex = $anchorIndex + 2

    $start = $d.Paragraphs($startParaIndex).Range.Start
    $end = $d.Paragraphs($endParaIndex).Range.End

    $r = $d.Range($start, $end)
    $r.Delete()
}
